$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume data (and the B36:C40 reordering)
# Leading apostrophe forces text interpretation, matching the original
# inline-string cells (General number format) without Excel coercing
# numeric-looking text into actual numbers.

$ws.Range("D2").Value = "'28.054.02"
$ws.Range("E2").Value = "'  +2.14%  "

$ws.Range("D3").Value = "'1.807.07"
$ws.Range("E3").Value = "'  +0.34%  "

$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "'  +0.56%  "

$ws.Range("D5").Value = "'339.58"
$ws.Range("E5").Value = "'  +0.23%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "'  +0.34%  "

$ws.Range("D7").Value = "'0.3911"
$ws.Range("E7").Value = "'  +2.73%  "

$ws.Range("D8").Value = "'0.3480"
$ws.Range("E8").Value = "'  +0.42%  "

$ws.Range("D9").Value = "'48.30"
$ws.Range("E9").Value = "'  -0.61%  "

$ws.Range("D10").Value = "'1.191"
$ws.Range("E10").Value = "'  -1.26%  "

$ws.Range("D11").Value = "'0.07566"
$ws.Range("E11").Value = "'  +0.43%  "

$ws.Range("D12").Value = "'1.005"
$ws.Range("E12").Value = "'  +0.53%  "

$ws.Range("D13").Value = "'22.04"
$ws.Range("E13").Value = "'  -0.62%  "

$ws.Range("D14").Value = "'6.505"
$ws.Range("E14").Value = "'  +0.11%  "

$ws.Range("D15").Value = "'1.818.66"
$ws.Range("E15").Value = "'  +1.13%  "

$ws.Range("D16").Value = "'7.138"
$ws.Range("E16").Value = "'  +0.49%  "

$ws.Range("D17").Value = "'0.00001101"
$ws.Range("E17").Value = "'  -0.50%  "

$ws.Range("D18").Value = "'0.06709"
$ws.Range("E18").Value = "'  +0.64%  "

$ws.Range("D19").Value = "'85.06"
$ws.Range("E19").Value = "'  +0.06%  "

$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "'  +0.57%  "

$ws.Range("D21").Value = "'17.75"
$ws.Range("E21").Value = "'  +1.70%  "

$ws.Range("D22").Value = "'6.572"
$ws.Range("E22").Value = "'  +0.52%  "

$ws.Range("D23").Value = "'28.098.06"
$ws.Range("E23").Value = "'  +2.38%  "

$ws.Range("D24").Value = "'12.41"
$ws.Range("E24").Value = "'  -1.66%  "

$ws.Range("D25").Value = "'2.419"

$ws.Range("D26").Value = "'1.486"
$ws.Range("E26").Value = "'  -1.17%  "

$ws.Range("D27").Value = "'2.535"
$ws.Range("E27").Value = "'  -1.73%  "

$ws.Range("D28").Value = "'21.29"
$ws.Range("E28").Value = "'  -1.30%  "

$ws.Range("D29").Value = "'154.29"
$ws.Range("E29").Value = "'  +1.46%  "

$ws.Range("D30").Value = "'2.023.76"
$ws.Range("E30").Value = "'  +1.12%  "

$ws.Range("D31").Value = "'136.00"
$ws.Range("E31").Value = "'  +1.41%  "

$ws.Range("D32").Value = "'4.038"
$ws.Range("E32").Value = "'  -0.32%  "

$ws.Range("D33").Value = "'6.143"
$ws.Range("E33").Value = "'  -0.22%  "

$ws.Range("D34").Value = "'0.08724"
$ws.Range("E34").Value = "'  +0.31%  "

$ws.Range("D35").Value = "'13.00"
$ws.Range("E35").Value = "'  -2.74%  "

$ws.Range("B36").Value = "'InternetComputer(DFINITY)"
$ws.Range("C36").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'5.471"
$ws.Range("E36").Value = "'  -0.16%  "

$ws.Range("B37").Value = "'VeChain"
$ws.Range("C37").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02427"
$ws.Range("E37").Value = "'  +3.40%  "

$ws.Range("B38").Value = "'Hedera"
$ws.Range("C38").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06549"
$ws.Range("E38").Value = "'  +2.15%  "

$ws.Range("B39").Value = "'TheSandbox"
$ws.Range("C39").Value = "'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6926"
$ws.Range("E39").Value = "'  -0.12%  "

$ws.Range("B40").Value = "'WEMIXTOKEN"
$ws.Range("C40").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "'1.611"
$ws.Range("E40").Value = "'  -1.97%  "

$ws.Range("D41").Value = "'0.2218"
$ws.Range("E41").Value = "'  +0.32%  "

$ws.Range("D42").Value = "'1.260"
$ws.Range("E42").Value = "'  -1.51%  "

$ws.Range("D43").Value = "'8.507"
$ws.Range("E43").Value = "'  -4.92%  "

$ws.Range("D44").Value = "'14.66"
$ws.Range("E44").Value = "'  +1.31%  "

$ws.Range("D45").Value = "'0.6445"
$ws.Range("E45").Value = "'  -0.55%  "

$ws.Range("D46").Value = "'3.873"
$ws.Range("E46").Value = "'  +0.16%  "

$ws.Range("D47").Value = "'2.154"
$ws.Range("E47").Value = "'  +0.17%  "

$ws.Range("D48").Value = "'131.06"
$ws.Range("E48").Value = "'  +0.26%  "

$ws.Range("D49").Value = "'0.07196"
$ws.Range("E49").Value = "'  -0.17%  "

$ws.Range("D50").Value = "'79.98"
$ws.Range("E50").Value = "'  -0.12%  "

$ws.Range("D51").Value = "'1.247"
$ws.Range("E51").Value = "'  +1.98%  "
